$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Retitle the English indicator name in C1 (old "Prevalence of anaemia..."
#    text is replaced by the new "Proportion of women with anemia..." text).
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "2.2.3 Proportion of women with anemia to the total population"

# ---------------------------------------------------------------------------
# 2) Add the 2022 data column (K) alongside the existing years (D:J).
#    Copy+PasteSpecial(formats) reproduces the number format / alignment /
#    border of the matching column-J cell on each row, then the value is
#    written on top of it.
# ---------------------------------------------------------------------------
function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# Header row: year 2022
Copy-Format "J3" "K3"
$ws.Range("K3").Value = 2022

# Row 4 - overall indicator (bold, one decimal)
Copy-Format "J4" "K4"
$ws.Range("K4").Value = 1.7

# Row 5 - "including:" sub-header, no value
Copy-Format "J5" "K5"

# Row 6 - children 0-14
Copy-Format "J6" "K6"
$ws.Range("K6").Value = 2.2000000000000002

# Row 7 - 15 years and older
Copy-Format "J7" "K7"
$ws.Range("K7").Value = 1.4

# Row 8 is a section header; it never had a K cell and the stray empty J8
# cell from the old layout is dropped entirely (see step 4 below).

# Row 9 - Kyrgyz Republic (bold-ish series like row 4)
Copy-Format "J9" "K9"
$ws.Range("K9").Value = 41.1

# Row 10 - Batken oblast
Copy-Format "J10" "K10"
$ws.Range("K10").Value = 65.90209110066462

# Row 11 - Djalal-Abad oblast
Copy-Format "J11" "K11"
$ws.Range("K11").Value = 55.941036331149498

# Row 12 - Ysyk-Kul oblast
Copy-Format "J12" "K12"
$ws.Range("K12").Value = 21.263715474839199

# Row 13 - Naryn oblast
Copy-Format "J13" "K13"
$ws.Range("K13").Value = 11.351981351981353

# Row 14 - Osh oblast
Copy-Format "J14" "K14"
$ws.Range("K14").Value = 32.279274699203526

# Row 15 - Talas oblast
Copy-Format "J15" "K15"
$ws.Range("K15").Value = 36.890901250539024

# Row 16 - Chui oblast
Copy-Format "J16" "K16"
$ws.Range("K16").Value = 32.421298573536646

# Row 17 - Bishkek city
Copy-Format "J17" "K17"
$ws.Range("K17").Value = 43.227712137486577

# Row 18 - Osh city (bottom, thick-bottom-bordered row)
Copy-Format "J5" "K18"
$ws.Range("K18").Value = 38.737482570668021
$k18Bottom = $ws.Range("K18").Borders.Item(9)
$k18Bottom.LineStyle = 1
$k18Bottom.Weight = -4138
$k18Bottom.Color = 0

# ---------------------------------------------------------------------------
# 3) Update the active selection to match the author's saved cursor position.
# ---------------------------------------------------------------------------
$ws.Range("M15").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4) Row 8 no longer carries a (stray, empty) J8 cell.
# ---------------------------------------------------------------------------
$ws.Range("J8").Clear() | Out-Null
